$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header "Produit" -> "Nom" ---
$ws.Cells.Item(1, 2).Value = "Nom"

# --- Rename "Légumes" -> "Légume" in the Type column (A2:A21) ---
for ($r = 2; $r -le 21; $r++) {
    $v = $ws.Cells.Item($r, 1).Value()
    if ($v -eq "Légumes") {
        $ws.Cells.Item($r, 1).Value = "Légume"
    }
}

# --- Add new (empty) row 22 with a bold E22 cell ---
$ws.Cells.Item(22, 5).Font.Bold = $true

# --- Apply a 2-decimal number format to Stock (C) and Prix /u (E) columns ---
$ws.Range("C2:C21").NumberFormat = "0.00"
$ws.Range("E2:E21").NumberFormat = "0.00"

# --- Update the active selection ---
[void]$ws.Range("L16").Select()
